$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "usercredentials"
$ws.Range("A1").Value = "admin"
$ws.Range("B1").Value = "manager"
$ws.Range("B1").Select()
